$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-09-04 20:57:53"

function Set-TextCell($row, $col, $text) {
    # Force text typing for values that look numeric (e.g. "6753975", "74.50")
    # by using the classic leading-apostrophe convention, then strip any
    # number-format/style drift that introduces so the cell keeps the
    # workbook's default style.
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# NOTE: columns D, E, I, J, K, L are already blank ("") on rows 11, 12 and 13
# before this script runs, and stay blank after the edit too -- so they are
# intentionally left untouched below (re-writing "" would replace the
# existing blank inline-string cell with a removed/empty cell instead).

# --- Row 11 becomes the old row 12's product (Severin Standgrill) ---
Set-TextCell 11 1 "6753975"
$ws.Cells.Item(11, 2).Value  = "Severin Standgrill mit Grillplatte PG 8563"
$ws.Cells.Item(11, 3).Value  = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-standgrill-mit-grillplatte-pg-8563/p/6753975"
$ws.Cells.Item(11, 6).Value  = 0
$ws.Cells.Item(11, 7).Value  = "Severin"
Set-TextCell 11 8 "74.50"
$ws.Cells.Item(11, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Cells.Item(11, 14).Value = "Severin Standgrill mit Grillplatte PG 8563 50% Aktion 74.50 Schweizer Franken statt 149.00 Schweizer Franken"
$ws.Cells.Item(11, 15).Value = $newTimestamp

# --- Row 12 becomes the old row 13's product (Tefal Dampfbuegeleisen) ---
Set-TextCell 12 1 "6821480"
$ws.Cells.Item(12, 2).Value  = "Tefal Dampfbügeleisen Express Protect SV9202S0 (7.5 Bar, Durilium)"
$ws.Cells.Item(12, 3).Value  = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/tefal-dampfbuegeleisen-express-protect-sv9202s0-75-bar-durilium/p/6821480"
$ws.Cells.Item(12, 6).Value  = 0
$ws.Cells.Item(12, 7).Value  = "Tefal"
Set-TextCell 12 8 "124.50"
$ws.Cells.Item(12, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Cells.Item(12, 14).Value = "Tefal Dampfbügeleisen Express Protect SV9202S0 (7.5 Bar, Durilium) 50% Aktion 124.50 Schweizer Franken statt 249.00 Schweizer Franken"
$ws.Cells.Item(12, 15).Value = $newTimestamp

# --- Row 13 becomes the old row 36's product (Varta Premium Light F20) ---
Set-TextCell 13 1 "6805078"
$ws.Cells.Item(13, 2).Value  = "Varta Premium Light F20"
$ws.Cells.Item(13, 3).Value  = "/de/haushalt-tier/elektroartikel-batterien/beleuchtung/taschenlampen-lichter/varta-premium-light-f20/p/6805078"
$ws.Cells.Item(13, 6).Value  = 0
$ws.Cells.Item(13, 7).Value  = "Varta"
Set-TextCell 13 8 "39.95"
$ws.Cells.Item(13, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'beleuchtung', 'taschenlampen-lichter']"
$ws.Cells.Item(13, 14).Value = "Varta Premium Light F20 39.95 Schweizer Franken"
$ws.Cells.Item(13, 15).Value = $newTimestamp

# --- Rows 14 through 35 keep their content; only the timestamp changes ---
for ($r = 14; $r -le 35; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# --- Remove the now-obsolete last row (old row 36, whose data moved to row 13) ---
$ws.Rows.Item(36).Delete()

# --- Rows 2 through 10 also only get the timestamp refreshed ---
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}
